$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.436.69'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '2.370.71'
$ws.Range("E3").Value = '  +3.11%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''309.50'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = '''105.22'
$ws.Range("E6").Value = '  +4.28%  '
$ws.Range("D7").Value = '''0.515'
$ws.Range("E7").Value = '  -4.00%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").Value = '''36.20'
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").Value = '''53.41'
$ws.Range("E11").Value = '  +2.51%  '
$ws.Range("D12").Value = '''0.0811'
$ws.Range("E12").Value = '  -1.52%  '
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("E14").Value = '  -1.77%  '
$ws.Range("D15").Value = '2.739.40'
$ws.Range("E15").Value = '  +3.07%  '
$ws.Range("D16").Value = '''15.63'
$ws.Range("E16").Value = '  +4.21%  '
$ws.Range("D17").Value = '2.374.20'
$ws.Range("E17").Value = '  +2.91%  '
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '43.379.56'
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("E20").Value = '  -4.07%  '
$ws.Range("D21").Value = '''6.31'
$ws.Range("E21").Value = '  +3.87%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '''68.31'
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").Value = '''241.36'
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("E25").Value = '  +1.85%  '
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '''25.74'
$ws.Range("E28").Value = '  +4.49%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.31'
$ws.Range("E29").Value = '  +9.10%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '''36.82'
$ws.Range("E30").Value = '  -4.23%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '''9.55'
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '''161.68'
$ws.Range("E32").Value = '  -3.69%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''5.27'
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D34").Value = '''18.43'
$ws.Range("E34").Value = '  +4.01%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").Value = '''4.77'
$ws.Range("E36").Value = '  +13.39%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '''2.54'
$ws.Range("E37").Value = '  +6.11%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '''3.10'
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.0740'
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '''1.94'
$ws.Range("E40").Value = '  +5.82%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.106'
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '''0.114'
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.017.16'
$ws.Range("E43").Value = '  +2.27%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''19.71'
$ws.Range("E44").Value = '  +2.90%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '''0.0290'
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '''3.15'
$ws.Range("E46").Value = '  +3.98%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '''10.57'
$ws.Range("E47").Value = '  +7.54%  '
$ws.Range("B48").Value = 'MultiversX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D48").Value = '''58.10'
$ws.Range("E48").Value = '  +4.49%  '
$ws.Range("B49").Value = 'HuobiToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D49").Value = '''2.97'
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.584.69'
$ws.Range("E50").Value = '  +2.31%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '''1.56'
$ws.Range("E51").Value = '  +1.61%  '
